$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vms")
$ws.Activate() | Out-Null

# Row 2 - new dev IIS VM (replaces old VSL-PRO-KCI-001 / nut-dmz-03 entry)
$ws.Range("A2").Value2 = "LAN"
$ws.Range("B2").Value2 = "VSL-DEV-IIS-001"
$ws.Range("C2").Value2 = "VSL-DEV-IIS-001"
$ws.Range("D2").Value2 = "VSL-DEV-IIS-001"
$ws.Range("E2").Value2 = "vsl-dev-iis-001"
$ws.Range("F2").Value2 = "nutanix.dc3"
$ws.Range("G2").Value2 = "pe_lu651"
$ws.Range("H2").Value2 = "rhel8-dc3"
$ws.Range("J2").Value2 = 'var.ahv_651_network["VLAN-20-Legacy-Server"]'
$ws.Range("L2").Value2 = 'var.ahv_651_storage["NUT_AHV_DC3_RH_PGSQL"]'
$ws.Range("N2").Value2 = 16384
$ws.Range("O2").Value2 = 4
$ws.Range("R2").Value2 = "172.17.20.148"
$ws.Range("T2").Value2 = "172.17.20.1"
$ws.Range("U2").Value2 = "DEV_TEST"

# Row 3 - new dev IDB VM (replaces old VSL-PRO-KCI-002 / nut-dmz-04 entry)
$ws.Range("A3").Value2 = "LAN"
$ws.Range("B3").Value2 = "VSL-DEV-IDB-001"
$ws.Range("C3").Value2 = "VSL-DEV-IDB-001"
$ws.Range("D3").Value2 = "VSL-DEV-IDB-001"
$ws.Range("E3").Value2 = "vsl-dev-idb-001"
$ws.Range("F3").Value2 = "nutanix.dc3"
$ws.Range("G3").Value2 = "pe_lu651"
$ws.Range("H3").Value2 = "rhel8-dc3"
$ws.Range("J3").Value2 = 'var.ahv_651_network["VLAN-20-Legacy-Server"]'
$ws.Range("L3").Value2 = 'var.ahv_651_storage["NUT_AHV_DC3_RH_PGSQL"]'
$ws.Range("N3").Value2 = 8192
$ws.Range("O3").Value2 = 2
$ws.Range("R3").Value2 = "172.17.20.149"
$ws.Range("T3").Value2 = "172.17.20.1"
$ws.Range("U3").Value2 = "DEV_TEST"

# Match the author's final selection
$ws.Range("Q3").Select() | Out-Null
